$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.151936173439026
$ws.Range("B1").Value = 2.024236440658569
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.2129807472229
$ws.Range("E1").Value = 1.157730579376221
